$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

Set-TextValue $ws "D2" "27.006.03"
Set-TextValue $ws "E2" "  -2.32%  "
Set-TextValue $ws "D3" "1.861.71"
Set-TextValue $ws "E3" "  -1.84%  "
Set-TextValue $ws "E4" "  -0.10%  "
Set-TextValue $ws "D5" "305.87"
Set-TextValue $ws "E5" "  -1.76%  "
Set-TextValue $ws "D6" "1.000"
Set-TextValue $ws "E6" "  +0.04%  "
Set-TextValue $ws "D7" "0.5099"
Set-TextValue $ws "E7" "  +2.66%  "
Set-TextValue $ws "D8" "0.3737"
Set-TextValue $ws "E8" "  -0.49%  "
Set-TextValue $ws "D9" "0.07111"
Set-TextValue $ws "E9" "  -1.63%  "
Set-TextValue $ws "D10" "0.8878"
Set-TextValue $ws "E10" "  +0.19%  "
Set-TextValue $ws "E11" "  -1.78%  "
Set-TextValue $ws "D12" "0.07553"
Set-TextValue $ws "E12" "  -0.67%  "
Set-TextValue $ws "D13" "1.860.65"
Set-TextValue $ws "E13" "  -4.31%  "
Set-TextValue $ws "E14" "  -2.49%  "
Set-TextValue $ws "D15" "89.20"
Set-TextValue $ws "E15" "  -2.46%  "
Set-TextValue $ws "D16" "1.001"
Set-TextValue $ws "E16" "  -0.15%  "
Set-TextValue $ws "D17" "0.000008349"
Set-TextValue $ws "E17" "  -3.58%  "
Set-TextValue $ws "D18" "14.09"
Set-TextValue $ws "E18" "  -2.32%  "
Set-TextValue $ws "D19" "1.001"
Set-TextValue $ws "E19" "  +0.18%  "
Set-TextValue $ws "D20" "27.044.98"
Set-TextValue $ws "E20" "  -2.37%  "
Set-TextValue $ws "D21" "5.052"
Set-TextValue $ws "E21" "  -1.42%  "
Set-TextValue $ws "D22" "2.091.78"
Set-TextValue $ws "E22" "  -5.80%  "
Set-TextValue $ws "E23" "  -2.70%  "
Set-TextValue $ws "D24" "6.462"
Set-TextValue $ws "E24" "  -1.47%  "
Set-TextValue $ws "D25" "148.82"
Set-TextValue $ws "E25" "  -2.88%  "
Set-TextValue $ws "D26" "1.837"
Set-TextValue $ws "E26" "  -0.22%  "
Set-TextValue $ws "E27" "  -1.36%  "
Set-TextValue $ws "D28" "2.081"
Set-TextValue $ws "E28" "  -4.90%  "
Set-TextValue $ws "D29" "112.80"
Set-TextValue $ws "E29" "  -1.41%  "
Set-TextValue $ws "D30" "4.675"
Set-TextValue $ws "E30" "  -2.90%  "
Set-TextValue $ws "D31" "4.649"
Set-TextValue $ws "E31" "  -2.40%  "
Set-TextValue $ws "D32" "0.09050"
Set-TextValue $ws "E32" "  +1.72%  "
Set-TextValue $ws "D33" "0.05120"
Set-TextValue $ws "E33" "  -2.82%  "
Set-TextValue $ws "D34" "3.055"
Set-TextValue $ws "E34" "  -4.00%  "
Set-TextValue $ws "D35" "1.153"
Set-TextValue $ws "E35" "  -5.42%  "
Set-TextValue $ws "D36" "0.7283"
Set-TextValue $ws "E36" "  -5.88%  "
Set-TextValue $ws "D37" "0.02045"
Set-TextValue $ws "E37" "  -0.86%  "
Set-TextValue $ws "D38" "3.046"
Set-TextValue $ws "E38" "  -0.15%  "
Set-TextValue $ws "D39" "2.469"
Set-TextValue $ws "E39" "  -4.94%  "
Set-TextValue $ws "D40" "1.069"
Set-TextValue $ws "E40" "  -1.67%  "
Set-TextValue $ws "D41" "0.5332"
Set-TextValue $ws "E41" "  -2.55%  "
Set-TextValue $ws "D42" "6.582"
Set-TextValue $ws "D43" "116.47"
Set-TextValue $ws "E43" "  +3.10%  "
Set-TextValue $ws "D44" "8.312"
Set-TextValue $ws "E44" "  -1.14%  "
Set-TextValue $ws "E45" "  -2.18%  "
Set-TextValue $ws "E46" "  +0.16%  "
Set-TextValue $ws "D47" "0.4611"
Set-TextValue $ws "E47" "  -2.98%  "
Set-TextValue $ws "E48" "  -3.64%  "
Set-TextValue $ws "E49" "  -2.60%  "
Set-TextValue $ws "D50" "36.61"
Set-TextValue $ws "E50" "  -0.23%  "
Set-TextValue $ws "D51" "63.99"
Set-TextValue $ws "E51" "  -3.81%  "
